# This script updates the simulation-derived transition probability matrix on Sheet1
# to reflect additional simulated games (per commit message: "added more games,
# sped up simulate game logic, and drafted optimization logic"). The underlying
# game-count data is external to the workbook, so the refreshed probabilities are
# written directly into the affected cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.2413793103448276
$ws.Range("C2").Value = 0.4827586206896552
$ws.Range("P2").Value = 0.1896551724137931
$ws.Range("S2").Value = 0.08620689655172414

# Row 3
$ws.Range("J3").Value = 0.07692307692307693
$ws.Range("P3").Value = 0.7692307692307693
$ws.Range("S3").Value = 0.1538461538461539

# Row 4
$ws.Range("J4").Value = 0.2
$ws.Range("P4").Value = 0.8

# Row 6
$ws.Range("B6").Value = 0.02083333333333333
$ws.Range("J6").Value = 0.1875
$ws.Range("Q6").Value = 0.08333333333333333
$ws.Range("R6").Value = 0.125
$ws.Range("S6").Value = 0.5833333333333334

# Row 7
$ws.Range("B7").Value = 0.02173913043478261
$ws.Range("F7").Value = 0.04347826086956522
$ws.Range("J7").Value = 0.08695652173913043
$ws.Range("Q7").Value = 0.108695652173913
$ws.Range("R7").Value = 0.02173913043478261
$ws.Range("S7").Value = 0.717391304347826

# Row 8
$ws.Range("B8").Value = 0.03389830508474576
$ws.Range("D8").Value = 0.00847457627118644
$ws.Range("F8").Value = 0.04661016949152542
$ws.Range("J8").Value = 0.0635593220338983
$ws.Range("O8").Value = 0.00847457627118644
$ws.Range("Q8").Value = 0.1059322033898305
$ws.Range("R8").Value = 0.1186440677966102
$ws.Range("S8").Value = 0.614406779661017

# Row 9
$ws.Range("B9").Value = 0.01754385964912281
$ws.Range("D9").Value = 0.01754385964912281
$ws.Range("F9").Value = 0.07017543859649122
$ws.Range("S9").Value = 0.5789473684210527

# Row 10
$ws.Range("B10").Value = 0.1107382550335571
$ws.Range("D10").Value = 0.006711409395973154
$ws.Range("F10").Value = 0.06040268456375839
$ws.Range("J10").Value = 0.087248322147651
$ws.Range("O10").Value = 0.01006711409395973
$ws.Range("Q10").Value = 0.1610738255033557
$ws.Range("R10").Value = 0.0738255033557047
$ws.Range("S10").Value = 0.4899328859060403

# Row 11
$ws.Range("G11").Value = 0.1186440677966102
$ws.Range("J11").Value = 0.03389830508474576
$ws.Range("K11").Value = 0.1864406779661017
$ws.Range("L11").Value = 0.6271186440677966
$ws.Range("S11").Value = 0.03389830508474576

# Row 12
$ws.Range("G12").Value = 0.8918918918918919
$ws.Range("J12").Value = 0.05405405405405406
$ws.Range("S12").Value = 0.05405405405405406

# Row 13
$ws.Range("G13").Value = 0.5384615384615384
$ws.Range("J13").Value = 0.3846153846153846
$ws.Range("S13").Value = 0.07692307692307693

# Row 15
$ws.Range("H15").Value = 0.3695652173913043
$ws.Range("I15").Value = 0.04347826086956522
$ws.Range("J15").Value = 0.2391304347826087
$ws.Range("K15").Value = 0.04347826086956522
$ws.Range("O15").Value = 0.02173913043478261
$ws.Range("S15").Value = 0.2826086956521739

# Row 16
$ws.Range("F16").Value = 0.02857142857142857
$ws.Range("H16").Value = 0.2285714285714286
$ws.Range("I16").Value = 0.08571428571428572
$ws.Range("J16").Value = 0.4285714285714285
$ws.Range("O16").Value = 0.05714285714285714
$ws.Range("S16").Value = 0.1714285714285714

# Row 17
$ws.Range("F17").Value = 0.03333333333333333
$ws.Range("H17").Value = 0.2333333333333333
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.3777777777777778
$ws.Range("K17").Value = 0.06666666666666667
$ws.Range("O17").Value = 0.05555555555555555
$ws.Range("S17").Value = 0.1333333333333333

# Row 18
$ws.Range("H18").Value = 0.2380952380952381
$ws.Range("I18").Value = 0.126984126984127
$ws.Range("K18").Value = 0.06349206349206349
$ws.Range("M18").Value = 0.03174603174603174
$ws.Range("O18").Value = 0.01587301587301587
$ws.Range("S18").Value = 0.1904761904761905

# Row 19
$ws.Range("F19").Value = 0.007692307692307693
$ws.Range("H19").Value = 0.3461538461538461
$ws.Range("I19").Value = 0.0673076923076923
$ws.Range("J19").Value = 0.2884615384615384
$ws.Range("K19").Value = 0.07115384615384615
$ws.Range("M19").Value = 0.02115384615384616
$ws.Range("O19").Value = 0.04423076923076923
$ws.Range("S19").Value = 0.1538461538461539

